$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AG2").Value = 1
